$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.970.19"
$ws.Range("E2").Value = "  -3.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.197.06"
$ws.Range("E3").Value = "  -7.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.73"
$ws.Range("E5").Value = "  -4.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "82.09"
$ws.Range("E6").Value = "  -5.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("E7").Value = "  -3.81%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.464"
$ws.Range("E9").Value = "  -6.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0771"
$ws.Range("E10").Value = "  -8.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "28.84"
$ws.Range("E11").Value = "  -5.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.17"
$ws.Range("E12").Value = "  -10.36%  "

$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.535.73"
$ws.Range("E14").Value = "  -7.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.23"
$ws.Range("E15").Value = "  -4.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.02"
$ws.Range("E16").Value = "  -7.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.203.41"
$ws.Range("E17").Value = "  -6.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.710"
$ws.Range("E18").Value = "  -6.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "38.891.95"
$ws.Range("E19").Value = "  -3.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("E20").Value = "  -5.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.68"
$ws.Range("E21").Value = "  -7.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.52"
$ws.Range("E22").Value = "  -5.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.17"
$ws.Range("E23").Value = "  -5.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.60"
$ws.Range("E24").Value = "  -3.54%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -7.27%  "

$ws.Range("E27").Value = "  -2.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.39"
$ws.Range("E28").Value = "  -5.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.03"
$ws.Range("E30").Value = "  -2.86%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "149.03"
$ws.Range("E31").Value = "  -3.79%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.80"
$ws.Range("E32").Value = "  -7.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  -7.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0691"
$ws.Range("E35").Value = "  -4.99%  "

$ws.Range("E36").Value = "  -4.00%  "

$ws.Range("E37").Value = "  -3.98%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0956"
$ws.Range("E38").Value = "  -4.13%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.63"
$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.75"
$ws.Range("E40").Value = "  -8.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.61"
$ws.Range("E41").Value = "  -5.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.67"
$ws.Range("E42").Value = "  -3.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.903.47"
$ws.Range("E43").Value = "  -3.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0257"
$ws.Range("E44").Value = "  -4.66%  "

$ws.Range("E45").Value = "  -15.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.61"
$ws.Range("E46").Value = "  -3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("E47").Value = "  -6.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.71"
$ws.Range("E48").Value = "  -11.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.423.31"
$ws.Range("E49").Value = "  -7.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.86"
$ws.Range("E50").Value = "  -3.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.38"
$ws.Range("E51").Value = "  -7.10%  "
